$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet index 1): insert new row 31 for the new ticket type, shifting rows 31-43 down to 32-44 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Rows(31).Insert()

# Copy formatting from row 30 into the newly-inserted blank row 31 so col A keeps its bold/centered/bordered style
$ws1.Range("A30:I30").Copy()
$ws1.Range("A31:I31").PasteSpecial(-4122)

# Populate the new row 31 with the new event/ticket data
$ws1.Range("A31").Value = 30
$ws1.Range("B31").Value = "2024-03-24"
$ws1.Range("C31").Value = "杭州·AD02动漫展--亦之紫F、L句号内场票"
$ws1.Range("D31").Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws1.Range("E31").Value = "2024.03.24 12:00-03.24 16:00"
$ws1.Range("F31").Value = 18
$ws1.Range("G31").Value = 258
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=81836"
$ws1.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202402/ecrRfQce1707375167618.jpeg"

# Refresh "想去人数" (interest count) values that changed between scrapes, sheet "展览"
$ws1.Range("F2").Value = 185
$ws1.Range("F5").Value = 958
$ws1.Range("F6").Value = 5183
$ws1.Range("F7").Value = 431
$ws1.Range("F8").Value = 607
$ws1.Range("F9").Value = 899
$ws1.Range("F10").Value = 810
$ws1.Range("F11").Value = 70
$ws1.Range("F12").Value = 26
$ws1.Range("F13").Value = 554
$ws1.Range("F14").Value = 8
$ws1.Range("F17").Value = 1718
$ws1.Range("F18").Value = 1440
$ws1.Range("F19").Value = 801
$ws1.Range("F20").Value = 294
$ws1.Range("F21").Value = 182
$ws1.Range("F22").Value = 293
$ws1.Range("F23").Value = 497
$ws1.Range("F24").Value = 127
$ws1.Range("F27").Value = 524
$ws1.Range("F28").Value = 2418
$ws1.Range("F29").Value = 168
$ws1.Range("F30").Value = 91
$ws1.Range("F32").Value = 82
$ws1.Range("F34").Value = 241
$ws1.Range("F40").Value = 620
$ws1.Range("F42").Value = 43

# --- Sheet "全部类型" (merged listing, sheet index 4): same interest-count refresh, no row insert here ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 185
$ws4.Range("F5").Value = 958
$ws4.Range("F7").Value = 5183
$ws4.Range("F8").Value = 431
$ws4.Range("F9").Value = 607
$ws4.Range("F12").Value = 899
$ws4.Range("F13").Value = 810
$ws4.Range("F15").Value = 70
$ws4.Range("F16").Value = 26
$ws4.Range("F17").Value = 554
$ws4.Range("F18").Value = 8
$ws4.Range("F22").Value = 1718
$ws4.Range("F23").Value = 1440
$ws4.Range("F24").Value = 801
$ws4.Range("F25").Value = 294
$ws4.Range("F26").Value = 182
$ws4.Range("F27").Value = 293
$ws4.Range("F29").Value = 497
$ws4.Range("F30").Value = 127
$ws4.Range("F32").Value = 524
$ws4.Range("F33").Value = 2418
$ws4.Range("F34").Value = 168
$ws4.Range("F35").Value = 91
$ws4.Range("F36").Value = 82
$ws4.Range("F38").Value = 241
$ws4.Range("F43").Value = 620
